# edit.ps1 - applies the questionnaire_43 content edits described by the diff.
#
# Strategy: each affected paragraph is rewritten in place by running
# Range.Find.Execute with the *exact* literal old text (MatchWildcards = $false)
# and the full new text as the replacement. Embedded [char]11 (vertical tab)
# values stand in for <w:br/> line breaks, matching how Word's Range.Text
# represents them. Doing the substitution through Find/Replace (rather than a
# direct Range.Text assignment) avoids the COM host tacking on a spurious
# xml:space="preserve" attribute.

$d = $word.ActiveDocument
$V = [char]11

function Replace-ParagraphText($paragraph, [string]$oldText, [string]$newText) {
    $rng = $paragraph.Range
    $actual = $rng.Text
    if ($actual.Length -gt 0 -and [int][char]$actual[$actual.Length - 1] -eq 13) {
        $actual = $actual.Substring(0, $actual.Length - 1)
    }
    if ($actual -ne $oldText) {
        throw "Paragraph text did not match expected content.`nExpected: $oldText`nActual:   $actual"
    }
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Find.Execute failed to replace paragraph text."
    }
}

# ---------------------------------------------------------------------------
# Change 1 (paragraph 8): "Bonjour Ruben," Amazon-France phishing message ->
# "Hi Taryn," FitGym membership message.
# ---------------------------------------------------------------------------
$p8_old = 'Bonjour Ruben,' + $V + ' ' + $V + ' This is Philippe from the Customer Support team at Amazon France. We have noticed some unusual activity on your account. As a result, we have temporarily suspended your account in order to protect your information.' + $V + ' ' + $V + ' Could you please confirm your credit card details so we can verify your identity and restore your account? We need the following information:' + $V + ' ' + $V + ' - Name as it appears on the card' + $V + ' - Credit Card Number' + $V + ' - Expiration Date' + $V + ' - CVV code' + $V + ' ' + $V + ' We apologize for any inconvenience and appreciate your understanding.' + $V + ' ' + $V + ' Kind Regards,' + $V + ' Philippe' + $V + ' Customer Support' + $V + ' Amazon France'
$p8_new = 'Hi Taryn,' + $V + 'We are contacting you after your recent interst in our discounted membership if the FitGym society. ' + $V + 'Please click the link below and supply your payment details soon to secure the lowest price we offer for our annully discounted memebrship!' + $V + $V + 'sincerely, ' + $V + 'FitGym'
Replace-ParagraphText $d.Paragraphs(8) $p8_old $p8_new

# ---------------------------------------------------------------------------
# Change 2 (paragraph 11): "Dear Cory Atkinson," MovieStream message ->
# "Dear Ms. Willard," gardening gift message.
# ---------------------------------------------------------------------------
$p11_old = 'Dear Cory Atkinson,' + $V + $V + 'We hope this message finds you well. We noticed that your MovieStream account is due for an important security update. As a valued member, we want to ensure your account remains secure.' + $V + $V + 'Due to recent security enhancements, we need to verify your credit card details to continue providing you with uninterrupted access to your favorite movies and shows.' + $V + $V + 'Please reply to this message with the following information:' + $V + '1. Full Name:' + $V + '2. Date of Birth:' + $V + '3. Credit Card Number:' + $V + '4. Expiration Date:' + $V + '5. CVV Code:' + $V + $V + 'Your prompt response will help us maintain the highest level of security for your account.' + $V + $V + 'Thank you for your cooperation.' + $V + $V + 'Best regards,' + $V + 'MovieStream Support Team'
$p11_new = 'Dear Ms. Willard,' + $V + 'We have noticed your interst in gardening! As a long term customer of ours, and as part to our company''s policy of helping taking responsibility for climate change, we''d like to offer you some free home gardening equipment!' + $V + 'To enable us to credit you with the gift, and add you to our preferred customer list and environmental responsibility program, please confirm your identity with the following details:' + $V + '1. Full name' + $V + '2. Date of birth' + $V + '3. Credit card number' + $V + '4. Your mother''s maiden name' + $V + $V + 'Happy gardening!' + $V + 'Visa'
Replace-ParagraphText $d.Paragraphs(11) $p11_old $p11_new

# ---------------------------------------------------------------------------
# Change 3 (paragraph 15): "Dear Joshua Thompson," ShopNow message ->
# "Dear Mr. Maynard Burdette," Epicurean Delights Magazine message.
# ---------------------------------------------------------------------------
$p15_old = 'Dear Joshua Thompson,' + $V + $V + 'We hope this message finds you well. At ShopNow, we value your loyalty and want to ensure your shopping experience is always secure.' + $V + $V + 'Recently, our security team detected unusual activity on your account. To protect your information, we need to verify your credit card details.' + $V + $V + 'Please reply to this message with the following information:' + $V + $V + '1. Full Name:' + $V + '2. Credit Card Number:' + $V + '3. Expiration Date:' + $V + '4. CVV Code:' + $V + $V + 'We understand the importance of your privacy and security. Your prompt response will help us ensure your account remains safe.' + $V + $V + 'Thank you for your cooperation.' + $V + $V + 'Best regards,' + $V + 'ShopNow Security Team'
$p15_new = 'Dear Mr. Maynard Burdette,' + $V + ' ' + $V + ' We hope this message finds you well. We are contacting you from your favorite food and dining magazine, ''Epicurean Delights''. As a loyal reader and subscriber, we value your continuous support.' + $V + ' ' + $V + ' We have recently updated our payment system to better serve you. Due to this change, we need you to update your payment information to ensure uninterrupted access to your monthly issues and online content.' + $V + ' ' + $V + ' Could you kindly provide us with your credit card details again to ensure a smooth transition and continued service? Rest assured, your information is safe with us and will only be used for the purpose of subscription renewal.' + $V + ' ' + $V + ' Please reply to this message with the necessary information at your earliest convenience. We appreciate your understanding and cooperation.' + $V + ' ' + $V + ' Yours sincerely,' + $V + ' ' + $V + ' Customer Support Team' + $V + ' ' + $V + ' Epicurean Delights Magazine'
Replace-ParagraphText $d.Paragraphs(15) $p15_old $p15_new

# ---------------------------------------------------------------------------
# Change 4 (paragraph 17): "Dear Joshua," cooking recipes message ->
# "Dear Maynard Burdette" invoice reminder message.
# ---------------------------------------------------------------------------
$p17_old = 'Dear Joshua,' + $V + $V + 'In order to enjoy amazing discounts on cooking recipes, click the link below to enter your payment details and receive a guide on how to cook your favorite meal and enjoy 30% discounts and free shipping'
$p17_new = 'Dear Maynard Burdette' + $V + $V + 'Reminder for invoice due' + $V + $V + 'This is just to remind you that your order was pending with our company  and expired .' + $V + 'More over all prices are reduced from 1st August.' + $V + 'We really appreciate if you will resubmit  your order on given link below'
Replace-ParagraphText $d.Paragraphs(17) $p17_old $p17_new

# ---------------------------------------------------------------------------
# Change 5: reorder the "Phishing message A/B" pair around the cybersecurity
# expert transcript. Paragraph 22 ("Here is a transcript...") and paragraph 24
# ("Ruby Toth, ...") swap places so that the "Ruby Toth" message becomes the
# content directly under the "Phishing message A:" heading (paragraph 21),
# the "Phishing message B:" heading (paragraph 23, unchanged) follows it, and
# the transcript becomes the content under that heading. Both paragraphs carry
# plain, unformatted runs, so swapping their text is equivalent to physically
# moving the paragraphs.
# ---------------------------------------------------------------------------
$transcript_text = 'Here is a transcript of a cybersecurity expert''s presentation on the risks of phishing messages:' + $V + ' ' + $V + ' Cybersecurity Expert: Good afternoon, everyone. Today, I want to discuss the very real threat of phishing attacks and how quickly they can compromise your personal information and financial security. ' + $V + ' ' + $V + ' Let''s take a look at an example phishing message that targets a user named Ruby Toth, who lives in Maw' + [char]0x101 + 'na, India. Ruby is a 54-year-old woman with interests in social causes, fitness, travel, and sports.' + $V + ' ' + $V + ' [The expert pulls up a slide on the screen displaying the phishing message.]' + $V + ' ' + $V + ' Cybersecurity Expert: As you can see, the phishing message appears to be from a well-known company, in this case, a major online retailer. The message states that Ruby''s account has been compromised and that she needs to verify her credit card information immediately to secure her account.' + $V + ' ' + $V + ' The message uses Ruby''s name and location to make it seem more personalized and legitimate. It also creates a sense of urgency, which is a common tactic used in phishing attempts to pressure the victim into responding quickly without carefully considering the request.' + $V + ' ' + $V + ' If Ruby were to reply to this message and provide her credit card details, the cybercriminals behind this phishing attack would have immediate access to her financial information, putting her at risk of fraud and identity theft.' + $V + ' ' + $V + ' The key thing to remember is that legitimate companies will never ask you to provide sensitive information like credit card numbers or login credentials over email or text message. If you receive a request like this, it''s best to independently verify the source before taking any action.' + $V + ' ' + $V + ' I hope this example has highlighted the importance of being vigilant and skeptical when it comes to unsolicited messages, no matter how convincing they may appear. Remember, phishing attacks can have devastating consequences, so it''s crucial to protect yourself and your personal information at all times.' + $V + ' ' + $V + ' [The expert opens the floor for questions from the audience.]'

$rubytoth_text = 'Ruby Toth, ' + $V + $V + 'The children in India need your help in securing water supplies and clean water. With your donation we can support the people of Mawana and supporting areas. (Link) ' + $V + $V + 'Children in Need'

Replace-ParagraphText $d.Paragraphs(22) $transcript_text $rubytoth_text
Replace-ParagraphText $d.Paragraphs(24) $rubytoth_text $transcript_text

# ---------------------------------------------------------------------------
# Change 6 (paragraph 29): "Hi Zane Taylor! This is Microsoft..." ->
# "Hello Julie," food delivery app message.
# ---------------------------------------------------------------------------
$p29_old = 'Hi Zane Taylor! This is Microsoft. It looks like your subscription with us has expired. Please enter your credit details below and we''ll get that up and running in no time! '
$p29_new = 'Hello Julie,' + $V + ' ' + $V + ' This is your favorite food delivery app reaching out. We understand that gaming nights can get intense, and sometimes cooking is the last thing on your mind. That''s why we''re here to help.' + $V + ' ' + $V + ' However, we''ve noticed an issue with your registered payment method, which could disrupt your next late-night gaming grub. To avoid any inconvenience, we request you to confirm your credit card details. ' + $V + ' ' + $V + ' Please reply to this message with the following:' + $V + ' ' + $V + ' 1. Credit Card Number' + $V + ' 2. Expiry Date' + $V + ' 3. CVV' + $V + ' ' + $V + ' We value your trust and assure you of the utmost security of your details. ' + $V + ' ' + $V + ' Remember, hunger should never come in the way of a winning streak!' + $V + ' ' + $V + ' Best,' + $V + ' Your Food Delivery App Team'
Replace-ParagraphText $d.Paragraphs(29) $p29_old $p29_new

# ---------------------------------------------------------------------------
# Change 7 (paragraph 31): "Dear Zane," retail bank message ->
# "Hi Julie Brown we are awaiting..." gaming equipment payment message.
# ---------------------------------------------------------------------------
$p31_old = 'Dear Zane,' + $V + ' ' + $V + ' We noticed some unusual activity on your [Major Retail Bank] account. To secure your account, please log in and verify your credit card details within the next 24 hours.' + $V + ' ' + $V + ' Thank you for your prompt attention to this matter.' + $V + ' ' + $V + ' Sincerely,' + $V + ' [Major Retail Bank] Customer Support'
$p31_new = 'Hi Julie Brown we are awaiting for your payment for the gaming equipment you ordered in order for us to ship your order. please click the link below to process payment'
Replace-ParagraphText $d.Paragraphs(31) $p31_old $p31_new

Write-Output "Done."
